$wb = $excel.ActiveWorkbook

$url = "https://github.com/OpenLocalizationTestOrg/oltest/blob/710146d34e6e66df6e5b17457cc9f42698cd0998/e2e/236f1989-d7a8-4e4b-9869-e1217105f4ec.md"
$mdName = "236f1989-d7a8-4e4b-9869-e1217105f4ec.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $url, "", "", $mdName)
$wsZh.Range("J2").Value = "236f1989-d7a8-4e4b-9869-e1217105f4ec.cc18edd51dd83142b3f30b9c4daea148e9e3d718.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-12 15:12:00"
$wsZh.Columns.Item(3).ColumnWidth = 29.17
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $url, "", "", $mdName)
$wsDe.Range("J2").Value = "236f1989-d7a8-4e4b-9869-e1217105f4ec.cc18edd51dd83142b3f30b9c4daea148e9e3d718.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-12 15:12:19"
$wsDe.Columns.Item(3).ColumnWidth = 29.17
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17

# ---- Overview sheet ----
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("E2").Value = "Handed back: in sync with en-US"
$wsOv.Range("F2").Value = "Handed back: in sync with en-US"
$wsOv.Columns.Item(5).ColumnWidth = 29.17
$wsOv.Columns.Item(6).ColumnWidth = 29.17
